# Implement the "affiliate system" changes into the books workbook:
#  - Correct the placeholder author names in column D (rows 2 & 3)
#  - Add a new "Price" column (E) with a header styled like the other
#    header cells but with left/right borders only, and per-row prices
#  - Update the selected cell to D4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix placeholder author names (were "Acc" / "Abn") ---
$ws.Range("D2").Value = "Alice Johnson"
$ws.Range("D3").Value = "Bob Williams"

# --- New "Price" header cell (E1), styled like the bold header row ---
$e1 = $ws.Range("E1")
$e1.Value = "Price"
$e1.Font.Bold = $true
$e1.HorizontalAlignment = -4108   # xlCenter
$e1.VerticalAlignment = -4160     # xlTop
$e1.Borders.Item(7).LineStyle = 1   # xlEdgeLeft, xlContinuous
$e1.Borders.Item(7).Weight = 2      # xlThin
$e1.Borders.Item(10).LineStyle = 1  # xlEdgeRight, xlContinuous
$e1.Borders.Item(10).Weight = 2     # xlThin

# --- New Price values for each book row ---
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 7
$ws.Range("E5").Value = 8
$ws.Range("E6").Value = 15

# --- Update the active selection to D4 ---
$ws.Range("D4").Select() | Out-Null
